$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - EFT / HESAPTAN EFT - Sube
$ws.Range("E3").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H3").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 4 - HESAPTAN EFT - ATM
$ws.Range("E4").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H4").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 5 - HESAPTAN EFT - Mobil
$ws.Range("E5").Value = "30,46 TL - 60,94 TL - 609,43 TL"
$ws.Range("H5").Value = "30,46 TL - 60,94 TL - 609,43 TL"

# Row 6 - DUZENLI EFT
$ws.Range("E6").Value = "6,09 TL - 12,19 TL - 152,35 TL"
$ws.Range("H6").Value = "6,09 TL - 12,19 TL - 152,35 TL"

# Row 8 - HESAPTAN HAVALE - Sube
$ws.Range("E8").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H8").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 9 - HESAPTAN HAVALE - ATM
$ws.Range("E9").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H9").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 10 - HESAPTAN HAVALE - Mobil
$ws.Range("E10").Value = "15,23 TL - 30,47 TL - 304,71 TL"
$ws.Range("H10").Value = "15,23 TL - 30,47 TL - 304,72 TL"

# Row 11 - DUZENLI HAVALE
$ws.Range("E11").Value = "3,04 TL - 6,09 TL - 76,17 TL"
$ws.Range("H11").Value = "3,05 TL - 6,1 TL - 76,18 TL"

# Row 13 - GELEN SWIFT
$ws.Range("C13").Value = "Hesaba: Asgari 0 TL | Azami 0,94 TL"
$ws.Range("D13").Value = "Hesaba: Asgari 1 TL | Azami 300 TL"
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("H13").Value = "Hesaba: Asgari 1 TL | Azami 6,09 TL"
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# Row 14 - GIDEN SWIFT - Mobil
$ws.Range("E14").Value = "1.660 TL - 1.660 TL"
$ws.Range("H14").Value = "2.100 TL - 4.300 TL"
